$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.068.53"
$ws.Range("E2").Value = "  -0.77%  "
$ws.Range("D3").Value = "2.612.42"
$ws.Range("E3").Value = "  -1.32%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.86"
$ws.Range("E5").Value = "  -1.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "164.58"
$ws.Range("E6").Value = "  -2.25%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -2.69%  "
$ws.Range("D9").Value = "2.611.01"
$ws.Range("E9").Value = "  -1.40%  "
$ws.Range("E10").Value = "  -5.28%  "
$ws.Range("E12").Value = "  -0.76%  "
$ws.Range("E13").Value = "  -0.80%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.22"
$ws.Range("E14").Value = "  -3.15%  "
$ws.Range("E16").Value = "  -2.99%  "
$ws.Range("D17").Value = "67.036.52"
$ws.Range("E17").Value = "  -0.55%  "
$ws.Range("D18").Value = "2.580.06"
$ws.Range("E18").Value = "  -2.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.72"
$ws.Range("E19").Value = "  -1.60%  "
$ws.Range("E20").Value = "  -1.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "354.97"
$ws.Range("E21").Value = "  -2.54%  "
$ws.Range("E22").Value = "  -3.21%  "
$ws.Range("E23").Value = "  -3.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.44"
$ws.Range("E24").Value = "  -5.10%  "
$ws.Range("E26").Value = "  -5.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "69.17"
$ws.Range("E27").Value = "  -2.26%  "
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("D30").Value = "0.0₃0992"
$ws.Range("E30").Value = "  -3.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "542.39"
$ws.Range("E31").Value = "  -2.61%  "
$ws.Range("E32").Value = "  -2.26%  "
$ws.Range("E34").Value = "  -3.17%  "
$ws.Range("E35").Value = "  -0.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("E37").Value = "  -4.63%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "157.51"
$ws.Range("E38").Value = "  +0.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.88"
$ws.Range("E39").Value = "  -2.85%  "
$ws.Range("E40").Value = "  -2.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.24"
$ws.Range("E41").Value = "  +1.68%  "
$ws.Range("E42").Value = "  -1.97%  "
$ws.Range("E43").Value = "  -3.53%  "
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("E45").Value = "  -5.17%  "
$ws.Range("E46").Value = "  -1.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "150.86"
$ws.Range("E47").Value = "  -2.00%  "
$ws.Range("E48").Value = "  -4.04%  "
$ws.Range("E49").Value = "  -3.06%  "
$ws.Range("E50").Value = "  -2.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0769"
$ws.Range("E51").Value = "  -1.53%  "
